$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8977755904197693
$ws.Range("B1").Value = 1.791259407997131
$ws.Range("C1").Value = 4.308572292327881
$ws.Range("D1").Value = 3.486431360244751
$ws.Range("E1").Value = 1.501083374023438
